# Update faturamento_diario_lojas.xlsx with new sales figures for columns S, T, U (and R3/R6)
# and recalculate the row totals in column AG (sum of B:AF for each row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bibi Cell Mundi
$ws.Range("S2").Value = 13662.3
$ws.Range("T2").Value = 3320.9
$ws.Range("AG2").Value = 168926.96

# Row 3 - Bibi Cell Vieiralves
$ws.Range("R3").Value = 6480
$ws.Range("S3").Value = 2690
$ws.Range("T3").Value = 3365
$ws.Range("AG3").Value = 84457.00999999999

# Row 4 - Bibi Cell Manauara
$ws.Range("S4").Value = 2521
$ws.Range("T4").Value = 5190
$ws.Range("U4").Value = 2247
$ws.Range("AG4").Value = 61697.9

# Row 5 - Bibi Cell Ponta Negra
$ws.Range("S5").Value = 2947
$ws.Range("T5").Value = 2106.75
$ws.Range("U5").Value = 1486
$ws.Range("AG5").Value = 54811.77

# Row 6 - total
$ws.Range("R6").Value = 20643.07
$ws.Range("S6").Value = 21820.3
$ws.Range("T6").Value = 13982.65
$ws.Range("U6").Value = 3733
$ws.Range("AG6").Value = 369893.64
